$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 1.077172790289468
$ws.Range("D2").Value = 46247219530.28053
$ws.Range("I2").Value = 42693213333.75719
$ws.Range("L2").Value = 0

# Row 3
$ws.Range("B3").Value = 1.800815659480598
$ws.Range("D3").Value = 53502294527.56234
$ws.Range("I3").Value = 47560712944.47166
$ws.Range("L3").Value = 0

# Row 4
$ws.Range("B4").Value = 1.130337585097727
$ws.Range("D4").Value = 46385782023.66711
$ws.Range("I4").Value = 42656364777.13512
$ws.Range("L4").Value = 0

# Row 5
$ws.Range("B5").Value = 1.794814984367979
$ws.Range("D5").Value = 53426790567.65892
$ws.Range("I5").Value = 47505007513.38911
$ws.Range("L5").Value = 0

# Row 6
$ws.Range("B6").Value = 1.077172790289468
$ws.Range("D6").Value = 46247219530.28053
$ws.Range("I6").Value = 42693213333.75719
$ws.Range("L6").Value = 0

# Row 7
$ws.Range("B7").Value = 1.800815659480598
$ws.Range("D7").Value = 53502294527.56234
$ws.Range("I7").Value = 47560712944.47166
$ws.Range("L7").Value = 0

# Row 8
$ws.Range("B8").Value = 1.130337585097727
$ws.Range("D8").Value = 46385782023.66711
$ws.Range("I8").Value = 42656364777.13512
$ws.Range("L8").Value = 0

# Row 9
$ws.Range("B9").Value = 1.794814984367979
$ws.Range("D9").Value = 53426790567.65892
$ws.Range("I9").Value = 47505007513.38911
$ws.Range("L9").Value = 0

# Row 10
$ws.Range("B10").Value = 0.6393717843845849
$ws.Range("D10").Value = 43196586691.04309
$ws.Range("I10").Value = 33951850879.80246
$ws.Range("L10").Value = 0

# Row 11
$ws.Range("B11").Value = 0.7552459685402271
$ws.Range("D11").Value = 48153337817.38026
$ws.Range("I11").Value = 37233166452.01514
$ws.Range("L11").Value = 0

# Row 12
$ws.Range("B12").Value = 0.6164307654473636
$ws.Range("D12").Value = 42578077213.51637
$ws.Range("I12").Value = 33665047718.53025
$ws.Range("L12").Value = 0

# Row 13
$ws.Range("B13").Value = 0.7574598594689712
$ws.Range("D13").Value = 47985365362.29326
$ws.Range("I13").Value = 37033183143.0697
$ws.Range("L13").Value = 0

# Row 14
$ws.Range("B14").Value = 0.6393717843845849
$ws.Range("D14").Value = 43196586691.04309
$ws.Range("I14").Value = 33951850879.80246
$ws.Range("L14").Value = 0

# Row 15
$ws.Range("B15").Value = 0.7552459685402206
$ws.Range("D15").Value = 48153337817.38026
$ws.Range("I15").Value = 37233166452.01523
$ws.Range("L15").Value = 0

# Row 16
$ws.Range("B16").Value = 0.6164307654473636
$ws.Range("D16").Value = 42578077213.51637
$ws.Range("I16").Value = 33665047718.53025
$ws.Range("L16").Value = 0

# Row 17
$ws.Range("B17").Value = 0.757459859468971
$ws.Range("D17").Value = 47985365362.29326
$ws.Range("I17").Value = 37033183143.0697
$ws.Range("L17").Value = 0
